$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 357 (shifts rows 357..437 down to 358..438)
$ws.Rows.Item(357).Insert()

# Populate the newly inserted row 357 with the new data.
$ws.Cells.Item(357, 1).Value = 8
$ws.Cells.Item(357, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(357, 3).Value = "Coquimbo"
$ws.Cells.Item(357, 4).Value = 44551
$ws.Cells.Item(357, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(357, 5).Value = 4
$ws.Cells.Item(357, 6).Value = 100112024
$ws.Cells.Item(357, 7).Value = "Choclo"
$ws.Cells.Item(357, 8).Value = "Dulce o Americano"
$ws.Cells.Item(357, 9).Value = "Primera"
$ws.Cells.Item(357, 10).Value = 28000
$ws.Cells.Item(357, 11).Value = 250
$ws.Cells.Item(357, 12).Value = 300
$ws.Cells.Item(357, 13).Value = 275
$ws.Cells.Item(357, 14).Value = "$/unidad"
$ws.Cells.Item(357, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(357, 16).Value = 275
$ws.Cells.Item(357, 17).Value = 1
$ws.Cells.Item(357, 18).Value = "Hortaliza"
